# staff navbar sync with the staff dashboard
#
# The receptionists sheet had a couple of columns (phone/age/salary) typed
# inconsistently between rows - row 2 (Qasim) stored them as numbers while
# row 3 (Maryam) stored them as text. Normalise row 2 to text (matching the
# rest of the sheet / the dashboard's expectations) and row 3 back to
# numbers, and refresh Qasim's privileges so "attendance" becomes the more
# specific "staff_attendance" (matching the staff dashboard's navbar).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qasim): phone/age/salary become text.
# A leading apostrophe is how Excel enters a numeric-looking value as text.
$ws.Range("D2").Value = "'923432928333"
$ws.Range("G2").Value = "'20"
$ws.Range("I2").Value = "'250000"

# Row 2 privileges: "attendance" -> "staff_attendance".
$ws.Range("L2").Value = "members,staff_attendance,payments,packages"

# Row 3 (Maryam): phone/age/salary become real numbers.
$ws.Range("D3").Value = 923432928333
$ws.Range("G3").Value = 20
$ws.Range("I3").Value = 250000
